$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489) on sheet ALC
if ($ws.Range("G2").Value2 -ne 5489) { Write-Host "WARNING: ALC G2 expected 5489 but found" $ws.Range("G2").Value2 }
$ws.Range("H2").Value = 588.3333
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 382.5
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 382.5
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -608.5

# Row 6 (Leve Item ID 4564) on sheet ALC
if ($ws.Range("G6").Value2 -ne 4564) { Write-Host "WARNING: ALC G6 expected 4564 but found" $ws.Range("G6").Value2 }
$ws.Range("H6").Value = 815.9459000000001
$ws.Range("I6").Value = 107.37037
$ws.Range("J6").Value = 2729.1
$ws.Range("K6").Value = 322.11111
$ws.Range("L6").Value = 8187.299999999999
$ws.Range("M6").Value = -210.11111
$ws.Range("N6").Value = -8411.299999999999

# Row 29 (Leve Item ID 4575) on sheet ALC
if ($ws.Range("G29").Value2 -ne 4575) { Write-Host "WARNING: ALC G29 expected 4575 but found" $ws.Range("G29").Value2 }
$ws.Range("H29").Value = 633.3333
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 1500
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 4500
$ws.Range("M29").Value = -319
$ws.Range("N29").Value = -5062

# Row 38 (Leve Item ID 4599) on sheet ALC
if ($ws.Range("G38").Value2 -ne 4599) { Write-Host "WARNING: ALC G38 expected 4599 but found" $ws.Range("G38").Value2 }
$ws.Range("H38").Value = 392.53845
$ws.Range("I38").Value = 300.27274
$ws.Range("J38").Value = 900
$ws.Range("K38").Value = 900.81822
$ws.Range("L38").Value = 2700
$ws.Range("M38").Value = -528.81822
$ws.Range("N38").Value = -3444

# Row 58 (Leve Item ID 4606) on sheet ALC
if ($ws.Range("G58").Value2 -ne 4606) { Write-Host "WARNING: ALC G58 expected 4606 but found" $ws.Range("G58").Value2 }
$ws.Range("H58").Value = 1804.9166
$ws.Range("I58").Value = 1365
$ws.Range("J58").Value = 2119.1428
$ws.Range("K58").Value = 4095
$ws.Range("L58").Value = 6357.428400000001
$ws.Range("M58").Value = -3945
$ws.Range("N58").Value = -6657.428400000001

# Row 87 (Leve Item ID 10651) on sheet ALC
if ($ws.Range("G87").Value2 -ne 10651) { Write-Host "WARNING: ALC G87 expected 10651 but found" $ws.Range("G87").Value2 }
$ws.Range("H87").Value = 23608.041
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 23608.041
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 23608.041
$ws.Range("N87").Value = -26104.041

# Row 90 (Leve Item ID 10651) on sheet ALC
if ($ws.Range("G90").Value2 -ne 10651) { Write-Host "WARNING: ALC G90 expected 10651 but found" $ws.Range("G90").Value2 }
$ws.Range("H90").Value = 23608.041
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 23608.041
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 70824.12300000001
$ws.Range("N90").Value = -83304.12300000001

# Row 137 (Leve Item ID 44013) on sheet ALC
if ($ws.Range("G137").Value2 -ne 44013) { Write-Host "WARNING: ALC G137 expected 44013 but found" $ws.Range("G137").Value2 }
$ws.Range("H137").Value = 1460.5
$ws.Range("I137").Value = 820.4
$ws.Range("J137").Value = 2100.6
$ws.Range("K137").Value = 2461.2
$ws.Range("L137").Value = 6301.799999999999
$ws.Range("M137").Value = 88.80000000000018
$ws.Range("N137").Value = -11401.8

# Row 138 (Leve Item ID 44169) on sheet ALC
if ($ws.Range("G138").Value2 -ne 44169) { Write-Host "WARNING: ALC G138 expected 44169 but found" $ws.Range("G138").Value2 }
$ws.Range("H138").Value = 2291.672
$ws.Range("I138").Value = 2105.2
$ws.Range("J138").Value = 2382.634
$ws.Range("K138").Value = 6315.599999999999
$ws.Range("L138").Value = 7147.902
$ws.Range("M138").Value = -1175.599999999999
$ws.Range("N138").Value = -17427.902


$ws = $wb.Worksheets.Item("ARM")
# Row 81 (Leve Item ID 10841) on sheet ARM
if ($ws.Range("G81").Value2 -ne 10841) { Write-Host "WARNING: ARM G81 expected 10841 but found" $ws.Range("G81").Value2 }
$ws.Range("H81").Value = 363333
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 363333
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 363333
$ws.Range("N81").Value = -365329

# Row 84 (Leve Item ID 10841) on sheet ARM
if ($ws.Range("G84").Value2 -ne 10841) { Write-Host "WARNING: ARM G84 expected 10841 but found" $ws.Range("G84").Value2 }
$ws.Range("H84").Value = 363333
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 363333
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 1089999
$ws.Range("N84").Value = -1099983


$ws = $wb.Worksheets.Item("BSM")
# Row 134 (Leve Item ID 43998) on sheet BSM
if ($ws.Range("G134").Value2 -ne 43998) { Write-Host "WARNING: BSM G134 expected 43998 but found" $ws.Range("G134").Value2 }
$ws.Range("H134").Value = 2191.2964
$ws.Range("I134").Value = 1754.0476
$ws.Range("J134").Value = 3721.6667
$ws.Range("K134").Value = 5262.142800000001
$ws.Range("L134").Value = 11165.0001
$ws.Range("M134").Value = -2727.142800000001
$ws.Range("N134").Value = -16235.0001


$ws = $wb.Worksheets.Item("CRP")
# Row 36 (Leve Item ID 1845) on sheet CRP
if ($ws.Range("G36").Value2 -ne 1845) { Write-Host "WARNING: CRP G36 expected 1845 but found" $ws.Range("G36").Value2 }
$ws.Range("H36").Value = 10000
$ws.Range("I36").Value = 10000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -9612

# Row 40 (Leve Item ID 1845) on sheet CRP
if ($ws.Range("G40").Value2 -ne 1845) { Write-Host "WARNING: CRP G40 expected 1845 but found" $ws.Range("G40").Value2 }
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 10000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -9840

# Row 108 (Leve Item ID 27087) on sheet CRP
if ($ws.Range("G108").Value2 -ne 27087) { Write-Host "WARNING: CRP G108 expected 27087 but found" $ws.Range("G108").Value2 }
$ws.Range("H108").Value = 28666.666
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 28666.666
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 28666.666
$ws.Range("N108").Value = -36346.666


$ws = $wb.Worksheets.Item("CUL")
# Row 39 (Leve Item ID 4712) on sheet CUL
if ($ws.Range("G39").Value2 -ne 4712) { Write-Host "WARNING: CUL G39 expected 4712 but found" $ws.Range("G39").Value2 }
$ws.Range("H39").Value = 2060
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2060
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 6180
$ws.Range("N39").Value = -6768

# Row 55 (Leve Item ID 4733) on sheet CUL
if ($ws.Range("G55").Value2 -ne 4733) { Write-Host "WARNING: CUL G55 expected 4733 but found" $ws.Range("G55").Value2 }
$ws.Range("H55").Value = 2058.8333
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 2058.8333
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 6176.499899999999
$ws.Range("N55").Value = -6530.499899999999

# Row 68 (Leve Item ID 12895) on sheet CUL
if ($ws.Range("G68").Value2 -ne 12895) { Write-Host "WARNING: CUL G68 expected 12895 but found" $ws.Range("G68").Value2 }
$ws.Range("H68").Value = 1407.5287
$ws.Range("I68").Value = 764.13043
$ws.Range("J68").Value = 1638.75
$ws.Range("K68").Value = 2292.39129
$ws.Range("L68").Value = 4916.25
$ws.Range("M68").Value = -1481.39129
$ws.Range("N68").Value = -6538.25

# Row 71 (Leve Item ID 12895) on sheet CUL
if ($ws.Range("G71").Value2 -ne 12895) { Write-Host "WARNING: CUL G71 expected 12895 but found" $ws.Range("G71").Value2 }
$ws.Range("H71").Value = 1407.5287
$ws.Range("I71").Value = 764.13043
$ws.Range("J71").Value = 1638.75
$ws.Range("K71").Value = 6877.173870000001
$ws.Range("L71").Value = 14748.75
$ws.Range("M71").Value = -2821.173870000001
$ws.Range("N71").Value = -22860.75

# Row 109 (Leve Item ID 27854) on sheet CUL
if ($ws.Range("G109").Value2 -ne 27854) { Write-Host "WARNING: CUL G109 expected 27854 but found" $ws.Range("G109").Value2 }
$ws.Range("H109").Value = 5787.6665
$ws.Range("I109").Value = 931.5
$ws.Range("J109").Value = 15500
$ws.Range("K109").Value = 2794.5
$ws.Range("L109").Value = 46500
$ws.Range("M109").Value = -1754.5
$ws.Range("N109").Value = -48580

# Row 112 (Leve Item ID 27855) on sheet CUL
if ($ws.Range("G112").Value2 -ne 27855) { Write-Host "WARNING: CUL G112 expected 27855 but found" $ws.Range("G112").Value2 }
$ws.Range("H112").Value = 7000
$ws.Range("I112").Value = 5500
$ws.Range("J112").Value = 8000
$ws.Range("K112").Value = 16500
$ws.Range("L112").Value = 24000
$ws.Range("M112").Value = -15392
$ws.Range("N112").Value = -26216

# Row 113 (Leve Item ID 27843) on sheet CUL
if ($ws.Range("G113").Value2 -ne 27843) { Write-Host "WARNING: CUL G113 expected 27843 but found" $ws.Range("G113").Value2 }
$ws.Range("H113").Value = 875.90247
$ws.Range("I113").Value = 509.15384
$ws.Range("J113").Value = 1511.6
$ws.Range("K113").Value = 1527.46152
$ws.Range("L113").Value = 4534.799999999999
$ws.Range("M113").Value = 642.5384799999999
$ws.Range("N113").Value = -8874.799999999999

# Row 123 (Leve Item ID 36037) on sheet CUL
if ($ws.Range("G123").Value2 -ne 36037) { Write-Host "WARNING: CUL G123 expected 36037 but found" $ws.Range("G123").Value2 }
$ws.Range("H123").Value = 10333.333
$ws.Range("I123").Value = 1000
$ws.Range("J123").Value = 15000
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 45000
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -49900

# Row 124 (Leve Item ID 36040) on sheet CUL
if ($ws.Range("G124").Value2 -ne 36040) { Write-Host "WARNING: CUL G124 expected 36040 but found" $ws.Range("G124").Value2 }
$ws.Range("H124").Value = 2129.3333
$ws.Range("I124").Value = 1833
$ws.Range("J124").Value = 2188.6
$ws.Range("K124").Value = 5499
$ws.Range("L124").Value = 6565.799999999999
$ws.Range("M124").Value = -589
$ws.Range("N124").Value = -16385.8

# Row 129 (Leve Item ID 36054) on sheet CUL
if ($ws.Range("G129").Value2 -ne 36054) { Write-Host "WARNING: CUL G129 expected 36054 but found" $ws.Range("G129").Value2 }
$ws.Range("H129").Value = 1228.5714
$ws.Range("I129").Value = 554.8
$ws.Range("J129").Value = 1841.091
$ws.Range("K129").Value = 1664.4
$ws.Range("L129").Value = 5523.272999999999
$ws.Range("M129").Value = 3335.6
$ws.Range("N129").Value = -15523.273

# Row 131 (Leve Item ID 36060) on sheet CUL
if ($ws.Range("G131").Value2 -ne 36060) { Write-Host "WARNING: CUL G131 expected 36060 but found" $ws.Range("G131").Value2 }
$ws.Range("H131").Value = 1178.5217
$ws.Range("I131").Value = 1008
$ws.Range("J131").Value = 1238.7059
$ws.Range("K131").Value = 3024
$ws.Range("L131").Value = 3716.1177
$ws.Range("M131").Value = 2016
$ws.Range("N131").Value = -13796.1177

# Row 133 (Leve Item ID 44073) on sheet CUL
if ($ws.Range("G133").Value2 -ne 44073) { Write-Host "WARNING: CUL G133 expected 44073 but found" $ws.Range("G133").Value2 }
$ws.Range("H133").Value = 16736.25
$ws.Range("I133").Value = 965
$ws.Range("J133").Value = 21993.334
$ws.Range("K133").Value = 2895
$ws.Range("L133").Value = 65980.00199999999
$ws.Range("M133").Value = 2165
$ws.Range("N133").Value = -76100.00199999999

# Row 134 (Leve Item ID 44074) on sheet CUL
if ($ws.Range("G134").Value2 -ne 44074) { Write-Host "WARNING: CUL G134 expected 44074 but found" $ws.Range("G134").Value2 }
$ws.Range("H134").Value = 3161.375
$ws.Range("I134").Value = 2275.5386
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 6826.6158
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -1756.6158
$ws.Range("N134").Value = -31140

# Row 137 (Leve Item ID 44088) on sheet CUL
if ($ws.Range("G137").Value2 -ne 44088) { Write-Host "WARNING: CUL G137 expected 44088 but found" $ws.Range("G137").Value2 }
$ws.Range("H137").Value = 6463
$ws.Range("I137").Value = 7175.684
$ws.Range("J137").Value = 4528.5713
$ws.Range("K137").Value = 21527.052
$ws.Range("L137").Value = 13585.7139
$ws.Range("M137").Value = -16427.052
$ws.Range("N137").Value = -23785.7139

# Row 138 (Leve Item ID 44105) on sheet CUL
if ($ws.Range("G138").Value2 -ne 44105) { Write-Host "WARNING: CUL G138 expected 44105 but found" $ws.Range("G138").Value2 }
$ws.Range("H138").Value = 1886.5555
$ws.Range("I138").Value = 761.6667
$ws.Range("J138").Value = 7511
$ws.Range("K138").Value = 2285.0001
$ws.Range("L138").Value = 22533
$ws.Range("M138").Value = 2854.9999
$ws.Range("N138").Value = -32813

# Row 139 (Leve Item ID 44102) on sheet CUL
if ($ws.Range("G139").Value2 -ne 44102) { Write-Host "WARNING: CUL G139 expected 44102 but found" $ws.Range("G139").Value2 }
$ws.Range("H139").Value = 2637.1428
$ws.Range("I139").Value = 2637.1428
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 7911.428400000001
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -2771.428400000001
$ws.Range("N139").ClearContents()

# Row 140 (Leve Item ID 44097) on sheet CUL
if ($ws.Range("G140").Value2 -ne 44097) { Write-Host "WARNING: CUL G140 expected 44097 but found" $ws.Range("G140").Value2 }
$ws.Range("H140").Value = 1912.6786
$ws.Range("I140").Value = 1325.0588
$ws.Range("J140").Value = 2820.818
$ws.Range("K140").Value = 3975.1764
$ws.Range("L140").Value = 8462.454000000002
$ws.Range("M140").Value = 1204.8236
$ws.Range("N140").Value = -18822.454


$ws = $wb.Worksheets.Item("GSM")
# Row 62 (Leve Item ID 11983) on sheet GSM
if ($ws.Range("G62").Value2 -ne 11983) { Write-Host "WARNING: GSM G62 expected 11983 but found" $ws.Range("G62").Value2 }
$ws.Range("H62").Value = 64992.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 64992.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 64992.5
$ws.Range("N62").Value = -66364.5

# Row 65 (Leve Item ID 11983) on sheet GSM
if ($ws.Range("G65").Value2 -ne 11983) { Write-Host "WARNING: GSM G65 expected 11983 but found" $ws.Range("G65").Value2 }
$ws.Range("H65").Value = 64992.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 64992.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 194977.5
$ws.Range("N65").Value = -201841.5

# Row 93 (Leve Item ID 18107) on sheet GSM
if ($ws.Range("G93").Value2 -ne 18107) { Write-Host "WARNING: GSM G93 expected 18107 but found" $ws.Range("G93").Value2 }
$ws.Range("H93").Value = 89947.5
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 89947.5
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 89947.5
$ws.Range("N93").Value = -93691.5

# Row 134 (Leve Item ID 42064) on sheet GSM
if ($ws.Range("G134").Value2 -ne 42064) { Write-Host "WARNING: GSM G134 expected 42064 but found" $ws.Range("G134").Value2 }
$ws.Range("H134").Value = 44000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 44000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 132000
$ws.Range("N134").Value = -137070


$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289) on sheet LTW
if ($ws.Range("G16").Value2 -ne 5289) { Write-Host "WARNING: LTW G16 expected 5289 but found" $ws.Range("G16").Value2 }
$ws.Range("H16").Value = 4083.111
$ws.Range("I16").Value = 3343.625
$ws.Range("J16").Value = 9999
$ws.Range("K16").Value = 3343.625
$ws.Range("L16").Value = 9999
$ws.Range("M16").Value = -3173.625
$ws.Range("N16").Value = -10339

# Row 133 (Leve Item ID 41903) on sheet LTW
if ($ws.Range("G133").Value2 -ne 41903) { Write-Host "WARNING: LTW G133 expected 41903 but found" $ws.Range("G133").Value2 }
$ws.Range("H133").Value = 40300
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 40300
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 40300
$ws.Range("N133").Value = -45360

# Row 135 (Leve Item ID 42036) on sheet LTW
if ($ws.Range("G135").Value2 -ne 42036) { Write-Host "WARNING: LTW G135 expected 42036 but found" $ws.Range("G135").Value2 }
$ws.Range("H135").Value = 55214.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 55214.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 55214.5
$ws.Range("N135").Value = -65354.5


$ws = $wb.Worksheets.Item("WVR")
# Row 123 (Leve Item ID 34127) on sheet WVR
if ($ws.Range("G123").Value2 -ne 34127) { Write-Host "WARNING: WVR G123 expected 34127 but found" $ws.Range("G123").Value2 }
$ws.Range("H123").Value = 24779.357
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24779.357
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24779.357
$ws.Range("N123").Value = -34579.357

